$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.952.51"
$ws.Range("E2").Value = "  +1.45%  "

$ws.Range("D3").Value = "1.675.79"
$ws.Range("E3").Value = "  +0.72%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "331.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.59%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9986"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.07%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3654"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.94%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "47.20"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.12%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3241"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.58%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.145"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.62%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07137"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.44%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9987"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.15%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.092"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.13%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.68"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.29%  "

$ws.Range("D15").Value = "1.668.92"
$ws.Range("E15").Value = "  +0.58%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.654"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.42%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001051"
$ws.Range("D17").Style = "Normal"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06549"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.05%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9985"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.20%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "78.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.86%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.88"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.00%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.915"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.10%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.84"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.85%  "

$ws.Range("D24").Value = "24.951.72"
$ws.Range("E24").Value = "  +1.61%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.445"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.79%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.403"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.01%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "148.44"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.16%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.71"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.13%  "

$ws.Range("D29").Value = "1.851.42"
$ws.Range("E29").Value = "  +0.34%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.04"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.52%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.185"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.32%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.088"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.73%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.804"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.77%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08486"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.41%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.659"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.15%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.31"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.81%  "

$ws.Range("E37").Value = "  -0.80%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06067"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.01%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02236"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.00%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.230"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.00%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2090"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.43%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.252"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.21%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9980"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.05%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5971"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.25%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.69"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +8.02%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.849"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.99%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5733"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.42%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.52"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.74%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.967"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.24%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07009"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.93%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.198"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.12%  "

